# Auto-generated edit script: updates market price/profit columns (H-N)
# on several worksheets, per the scheduled-runner diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 157
$ws.Range("I12").Value = 145.33333
$ws.Range("J12").Value = 192
$ws.Range("K12").Value = 145.33333
$ws.Range("L12").Value = 192
$ws.Range("M12").Value = 24.66667000000001
$ws.Range("N12").Value = -532
$ws.Range("H33").Value = 889.4545000000001
$ws.Range("I33").Value = 958.4
$ws.Range("K33").Value = 958.4
$ws.Range("M33").Value = -729.4
$ws.Range("H43").Value = 999.5
$ws.Range("J43").Value = 1000
$ws.Range("L43").Value = 1000
$ws.Range("N43").Value = -1138
$ws.Range("H80").Value = 3253.35
$ws.Range("I80").Value = 1569.6364
$ws.Range("J80").Value = 5311.222
$ws.Range("K80").Value = 4708.9092
$ws.Range("L80").Value = 15933.666
$ws.Range("M80").Value = -3710.9092
$ws.Range("N80").Value = -17929.666
$ws.Range("H83").Value = 3253.35
$ws.Range("I83").Value = 1569.6364
$ws.Range("J83").Value = 5311.222
$ws.Range("K83").Value = 14126.7276
$ws.Range("L83").Value = 47800.998
$ws.Range("M83").Value = -9134.7276
$ws.Range("N83").Value = -57784.998
$ws.Range("H92").Value = 480.9091
$ws.Range("I92").Value = 414.14285
$ws.Range("K92").Value = 414.14285
$ws.Range("M92").Value = 833.85715
$ws.Range("H99").Value = 523.6
$ws.Range("I99").Value = 289.5
$ws.Range("J99").Value = 679.6667
$ws.Range("K99").Value = 868.5
$ws.Range("L99").Value = 2039.0001
$ws.Range("M99").Value = 629.5
$ws.Range("N99").Value = -5035.0001
$ws.Range("H101").Value = 14286187
$ws.Range("I101").Value = 20000422
$ws.Range("J101").Value = 599.5
$ws.Range("K101").Value = 60001266
$ws.Range("L101").Value = 1798.5
$ws.Range("M101").Value = -59999644
$ws.Range("N101").Value = -5042.5
$ws.Range("H115").Value = 180
$ws.Range("I115").Value = 180
$ws.Range("K115").Value = 540
$ws.Range("M115").Value = 1027
$ws.Range("H116").Value = 5955.3335
$ws.Range("I116").Value = 1999.6666
$ws.Range("J116").Value = 7933.1665
$ws.Range("K116").Value = 1999.6666
$ws.Range("L116").Value = 7933.1665
$ws.Range("M116").Value = 1442.3334
$ws.Range("N116").Value = -14817.1665
$ws.Range("H138").Value = 2070.3044
$ws.Range("J138").Value = 1936.8
$ws.Range("L138").Value = 5810.4
$ws.Range("N138").Value = -16090.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 1869
$ws.Range("J21").Value = 1697
$ws.Range("L21").Value = 1697
$ws.Range("N21").Value = -2445
$ws.Range("H24").Value = 35116.332
$ws.Range("J24").Value = 35116.332
$ws.Range("L24").Value = 35116.332
$ws.Range("N24").Value = -35864.332
$ws.Range("H31").Value = 2000
$ws.Range("I31").Value = 2000
$ws.Range("K31").Value = 2000
$ws.Range("M31").Value = -1706
$ws.Range("H32").Value = 3293.7856
$ws.Range("I32").Value = 3045.4443
$ws.Range("K32").Value = 3045.4443
$ws.Range("M32").Value = -2758.4443
$ws.Range("H96").Value = 37495
$ws.Range("J96").Value = 37495
$ws.Range("L96").Value = 37495
$ws.Range("N96").Value = -42987
$ws.Range("H100").Value = 35116.332
$ws.Range("J100").Value = 35116.332
$ws.Range("L100").Value = 35116.332
$ws.Range("N100").Value = -37280.332
$ws.Range("H132").Value = 783.6111
$ws.Range("I132").Value = 783.6111
$ws.Range("K132").Value = 2350.8333
$ws.Range("M132").Value = 179.1667000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 3348082.8
$ws.Range("J4").Value = 26166
$ws.Range("L4").Value = 26166
$ws.Range("N4").Value = -26390
$ws.Range("H22").Value = 918.8125
$ws.Range("I22").Value = 605
$ws.Range("J22").Value = 1322.2858
$ws.Range("K22").Value = 605
$ws.Range("L22").Value = 1322.2858
$ws.Range("M22").Value = -255
$ws.Range("N22").Value = -2022.2858
$ws.Range("H62").Value = 3197.2856
$ws.Range("I62").Value = 1776.8
$ws.Range("K62").Value = 1776.8
$ws.Range("M62").Value = -1152.8
$ws.Range("H65").Value = 3197.2856
$ws.Range("I65").Value = 1776.8
$ws.Range("K65").Value = 8884
$ws.Range("M65").Value = -5764
$ws.Range("H88").Value = 22213.334
$ws.Range("J88").Value = 22213.334
$ws.Range("L88").Value = 22213.334
$ws.Range("N88").Value = -23025.334
$ws.Range("H91").Value = 22213.334
$ws.Range("J91").Value = 22213.334
$ws.Range("L91").Value = 22213.334
$ws.Range("N91").Value = -25021.334
$ws.Range("H134").Value = 945.0714
$ws.Range("I134").Value = 940.9231
$ws.Range("J134").Value = 999
$ws.Range("K134").Value = 2822.7693
$ws.Range("L134").Value = 2997
$ws.Range("M134").Value = -287.7692999999999
$ws.Range("N134").Value = -8067
$ws.Range("H141").Value = 140994.73
$ws.Range("J141").Value = 150595.3
$ws.Range("L141").Value = 150595.3
$ws.Range("N141").Value = -160955.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 786.0625
$ws.Range("I14").Value = 786.0625
$ws.Range("K14").Value = 2358.1875
$ws.Range("M14").Value = -2185.1875
$ws.Range("H23").Value = 321.625
$ws.Range("J23").Value = 729
$ws.Range("L23").Value = 2187
$ws.Range("N23").Value = -2657
$ws.Range("H38").Value = 39.4
$ws.Range("I38").Value = 27.11111
$ws.Range("J38").Value = 150
$ws.Range("K38").Value = 81.33333
$ws.Range("L38").Value = 450
$ws.Range("M38").Value = 265.66667
$ws.Range("N38").Value = -1144
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H109").Value = 1615.5
$ws.Range("I109").Value = 1615.5
$ws.Range("K109").Value = 4846.5
$ws.Range("M109").Value = -3806.5
$ws.Range("H117").Value = 16084.286
$ws.Range("I117").Value = 300
$ws.Range("K117").Value = 900
$ws.Range("M117").Value = 2542
$ws.Range("H140").Value = 7679.222
$ws.Range("I140").Value = 1196.7059
$ws.Range("K140").Value = 3590.1177
$ws.Range("M140").Value = 1589.8823

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1733.1666
$ws.Range("I80").Value = 1699.75
$ws.Range("J80").Value = 1800
$ws.Range("K80").Value = 1699.75
$ws.Range("L80").Value = 1800
$ws.Range("M80").Value = -701.75
$ws.Range("N80").Value = -3796
$ws.Range("H83").Value = 1733.1666
$ws.Range("I83").Value = 1699.75
$ws.Range("J83").Value = 1800
$ws.Range("K83").Value = 8498.75
$ws.Range("L83").Value = 9000
$ws.Range("M83").Value = -3506.75
$ws.Range("N83").Value = -18984
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 1880.6154
$ws.Range("I132").Value = 1880.6154
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5641.8462
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3111.8462
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1707.5
$ws.Range("I93").Value = 1538.6666
$ws.Range("J93").Value = 1902.3077
$ws.Range("K93").Value = 1538.6666
$ws.Range("L93").Value = 1902.3077
$ws.Range("M93").Value = -290.6666
$ws.Range("N93").Value = -4398.3077
$ws.Range("H132").Value = 2055.05
$ws.Range("I132").Value = 1727.8889
$ws.Range("J132").Value = 4999.5
$ws.Range("K132").Value = 5183.6667
$ws.Range("L132").Value = 14998.5
$ws.Range("M132").Value = -2653.6667
$ws.Range("N132").Value = -20058.5
$ws.Range("H136").Value = 3292.5557
$ws.Range("I136").Value = 2254.8333
$ws.Range("J136").Value = 5368
$ws.Range("K136").Value = 6764.499899999999
$ws.Range("L136").Value = 16104
$ws.Range("M136").Value = -4214.499899999999
$ws.Range("N136").Value = -21204

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2284.8
$ws.Range("I132").Value = 2363.9583
$ws.Range("J132").Value = 385
$ws.Range("K132").Value = 7091.874899999999
$ws.Range("L132").Value = 1155
$ws.Range("M132").Value = -4561.874899999999
$ws.Range("N132").Value = -6215

